$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C2:C11) from 45221 to 45224
$ws.Range("C2:C11").Value = 45224
